# Updated cryptos list - Mon May  1 07:00:07 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue "D2" "28.612.55"
$ws.Range("E2").Value = "  -3.18%  "

# Row 3 - Ethereum
Set-TextValue "D3" "1.850.97"
$ws.Range("E3").Value = "  -3.75%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.99%  "

# Row 5 - BNB
Set-TextValue "D5" "334.16"
$ws.Range("E5").Value = "  +2.49%  "

# Row 6 - USDC
$ws.Range("E6").Value = "  -0.85%  "

# Row 7
Set-TextValue "D7" "0.4644"
$ws.Range("E7").Value = "  -3.60%  "

# Row 8
Set-TextValue "D8" "0.3918"
$ws.Range("E8").Value = "  -3.79%  "

# Row 9
Set-TextValue "D9" "46.50"
$ws.Range("E9").Value = "  -2.58%  "

# Row 10
Set-TextValue "D10" "0.07917"
$ws.Range("E10").Value = "  -3.67%  "

# Row 11
Set-TextValue "D11" "0.9860"
$ws.Range("E11").Value = "  -2.42%  "

# Row 12
Set-TextValue "D12" "22.29"
$ws.Range("E12").Value = "  -5.80%  "

# Row 13
Set-TextValue "D13" "1.968.58"
$ws.Range("E13").Value = "  +3.95%  "

# Row 14
Set-TextValue "D14" "5.850"
$ws.Range("E14").Value = "  -3.80%  "

# Row 15
Set-TextValue "D15" "7.015"
$ws.Range("E15").Value = "  -3.52%  "

# Row 16
Set-TextValue "D16" "0.06865"
$ws.Range("E16").Value = "  -0.21%  "

# Row 17
Set-TextValue "D17" "87.85"
$ws.Range("E17").Value = "  -4.20%  "

# Row 18
$ws.Range("E18").Value = "  -0.95%  "

# Row 19
$ws.Range("E19").Value = "  -2.97%  "

# Row 20
Set-TextValue "D20" "17.14"
$ws.Range("E20").Value = "  -2.82%  "

# Row 21
$ws.Range("E21").Value = "  -0.80%  "

# Row 22
Set-TextValue "D22" "28.646.69"
$ws.Range("E22").Value = "  -3.10%  "

# Row 23
Set-TextValue "D23" "5.395"
$ws.Range("E23").Value = "  -5.15%  "

# Row 24
Set-TextValue "D24" "11.32"
$ws.Range("E24").Value = "  -5.05%  "

# Row 25
Set-TextValue "D25" "2.221.95"
$ws.Range("E25").Value = "  +4.54%  "

# Row 26
$ws.Range("E26").Value = "  -2.48%  "

# Row 27
Set-TextValue "D27" "153.07"
$ws.Range("E27").Value = "  -1.85%  "

# Row 28
Set-TextValue "D28" "19.42"
$ws.Range("E28").Value = "  -3.08%  "

# Row 29
Set-TextValue "D29" "6.109"
$ws.Range("E29").Value = "  -5.76%  "

# Row 30
Set-TextValue "D30" "2.018"
$ws.Range("E30").Value = "  -3.87%  "

# Row 31
Set-TextValue "D31" "117.53"
$ws.Range("E31").Value = "  -2.56%  "

# Row 32
Set-TextValue "D32" "0.9798"
$ws.Range("E32").Value = "  -3.67%  "

# Row 33
Set-TextValue "D33" "0.09415"
$ws.Range("E33").Value = "  -2.28%  "

# Row 34
Set-TextValue "D34" "5.369"
$ws.Range("E34").Value = "  -4.67%  "

# Row 35
Set-TextValue "D35" "3.478"
$ws.Range("E35").Value = "  -2.01%  "

# Row 36
Set-TextValue "D36" "1.351"
$ws.Range("E36").Value = "  -1.88%  "

# Row 37
Set-TextValue "D37" "0.06151"
$ws.Range("E37").Value = "  -3.52%  "

# Row 38
$ws.Range("E38").Value = "  -4.11%  "

# Row 39
Set-TextValue "D39" "1.163"
$ws.Range("E39").Value = "  -2.02%  "

# Row 40
Set-TextValue "D40" "0.5713"
$ws.Range("E40").Value = "  -4.01%  "

# Row 41
Set-TextValue "D41" "7.632"
$ws.Range("E41").Value = "  -3.33%  "

# Row 42
Set-TextValue "D42" "10.15"
$ws.Range("E42").Value = "  -5.92%  "

# Row 43
$ws.Range("E43").Value = "  -2.72%  "

# Row 44
Set-TextValue "D44" "2.392"
$ws.Range("E44").Value = "  -3.37%  "

# Row 45
$ws.Range("E45").Value = "  -2.41%  "

# Row 46 - swaps with row 47 (Decentraland moves up to rank 46)
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
Set-TextValue "D46" "0.5398"
$ws.Range("E46").Value = "  -3.12%  "

# Row 47 - EnergySwap moves down to rank 47
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D47" "11.80"
$ws.Range("E47").Value = "  -4.70%  "

# Row 48
Set-TextValue "D48" "0.07157"
$ws.Range("E48").Value = "  -4.57%  "

# Row 49
$ws.Range("E49").Value = "  -2.15%  "

# Row 50
Set-TextValue "D50" "114.19"
$ws.Range("E50").Value = "  -4.19%  "

# Row 51
Set-TextValue "D51" "42.81"
$ws.Range("E51").Value = "  +2.20%  "
